# This script applies the latest Universalis market-board price refresh
# to the per-job Leve profit tables (columns H-N) across all 8 job sheets,
# mirroring the workbook's scheduled data-refresh runner.
#   H = currentAveragePrice        K = LevePriceNQ (= I * Leve Amount)
#   I = currentAveragePriceNQ      L = LevePriceHQ (= J * Leve Amount)
#   J = currentAveragePriceHQ      M = LeveProfitNQ (= Leve Gil - K)
#                                   N = LeveProfitHQ (= -2*Leve Gil - L)

$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# Row 88: The Grave of Hemlock Groves
$ws.Cells.Item(88, 8).Value = 2555.5454
$ws.Cells.Item(88, 9).Value = 3001.5
$ws.Cells.Item(88, 10).Value = 2456.4443
$ws.Cells.Item(88, 11).Value = 3001.5
$ws.Cells.Item(88, 12).Value = 2456.4443
$ws.Cells.Item(88, 13).Value = -2595.5
$ws.Cells.Item(88, 14).Value = -3268.4443
# Row 91: Dappling the Highlands (L)
$ws.Cells.Item(91, 8).Value = 2555.5454
$ws.Cells.Item(91, 9).Value = 3001.5
$ws.Cells.Item(91, 10).Value = 2456.4443
$ws.Cells.Item(91, 11).Value = 3001.5
$ws.Cells.Item(91, 12).Value = 2456.4443
$ws.Cells.Item(91, 13).Value = -1597.5
$ws.Cells.Item(91, 14).Value = -5264.4443
# Row 137: Cutting Edge of Culinary Quality
$ws.Cells.Item(137, 8).Value = 1234.1395
$ws.Cells.Item(137, 9).Value = 1238.1515
$ws.Cells.Item(137, 10).Value = 1220.9
$ws.Cells.Item(137, 11).Value = 3714.4545
$ws.Cells.Item(137, 12).Value = 3662.7
$ws.Cells.Item(137, 13).Value = -1164.4545
$ws.Cells.Item(137, 14).Value = -8762.700000000001
# Row 138: All-night Crafting
$ws.Cells.Item(138, 8).Value = 2177.5615
$ws.Cells.Item(138, 9).Value = 1147.375
$ws.Cells.Item(138, 11).Value = 3442.125
$ws.Cells.Item(138, 13).Value = 1697.875

$ws = $wb.Worksheets.Item("ARM")
# Row 26: Night Squawker
$ws.Cells.Item(26, 8).Value = 3602
$ws.Cells.Item(26, 9).Value = 1802.25
$ws.Cells.Item(26, 10).Value = 18000
$ws.Cells.Item(26, 11).Value = 1802.25
$ws.Cells.Item(26, 12).Value = 18000
$ws.Cells.Item(26, 13).Value = -1472.25
$ws.Cells.Item(26, 14).Value = -18660
# Row 61: Dealing with the Tough Stuff
$ws.Cells.Item(61, 8).Value = 2168.24
$ws.Cells.Item(61, 9).Value = 1860.3
$ws.Cells.Item(61, 11).Value = 1860.3
$ws.Cells.Item(61, 13).Value = -1648.3
# Row 63: Rivets Run through It
$ws.Cells.Item(63, 8).Value = 3109
$ws.Cells.Item(63, 9).Value = 2530.8
$ws.Cells.Item(63, 11).Value = 2530.8
$ws.Cells.Item(63, 13).Value = -1844.8
# Row 66: A Riveting Revival (L)
$ws.Cells.Item(66, 8).Value = 3109
$ws.Cells.Item(66, 9).Value = 2530.8
$ws.Cells.Item(66, 11).Value = 12654
$ws.Cells.Item(66, 13).Value = -9222
# Row 74: As the Bolt Flies
$ws.Cells.Item(74, 8).Value = 1500
$ws.Cells.Item(74, 9).Value = 833.3333
$ws.Cells.Item(74, 10).Value = 1900
$ws.Cells.Item(74, 11).Value = 833.3333
$ws.Cells.Item(74, 12).Value = 1900
$ws.Cells.Item(74, 13).Value = 40.66669999999999
$ws.Cells.Item(74, 14).Value = -3648
# Row 77: Heavy Metal Banned (L)
$ws.Cells.Item(77, 8).Value = 1500
$ws.Cells.Item(77, 9).Value = 833.3333
$ws.Cells.Item(77, 10).Value = 1900
$ws.Cells.Item(77, 11).Value = 4166.6665
$ws.Cells.Item(77, 12).Value = 9500
$ws.Cells.Item(77, 13).Value = 201.3334999999997
$ws.Cells.Item(77, 14).Value = -18236
# Row 88: The Mast Chance
$ws.Cells.Item(88, 8).Value = 2875.5
$ws.Cells.Item(88, 9).Value = 2367.6667
$ws.Cells.Item(88, 10).Value = 3383.3333
$ws.Cells.Item(88, 11).Value = 2367.6667
$ws.Cells.Item(88, 12).Value = 3383.3333
$ws.Cells.Item(88, 13).Value = -1961.6667
$ws.Cells.Item(88, 14).Value = -4195.3333
# Row 91: The Rose and the Riveter (L)
$ws.Cells.Item(91, 8).Value = 2875.5
$ws.Cells.Item(91, 9).Value = 2367.6667
$ws.Cells.Item(91, 10).Value = 3383.3333
$ws.Cells.Item(91, 11).Value = 2367.6667
$ws.Cells.Item(91, 12).Value = 3383.3333
$ws.Cells.Item(91, 13).Value = -963.6667000000002
$ws.Cells.Item(91, 14).Value = -6191.3333
# Row 132: Don't Bore Me, Ore Me
$ws.Cells.Item(132, 8).Value = 4118.646
$ws.Cells.Item(132, 9).Value = 5098.933
$ws.Cells.Item(132, 10).Value = 2484.8333
$ws.Cells.Item(132, 11).Value = 15296.799
$ws.Cells.Item(132, 12).Value = 7454.499899999999
$ws.Cells.Item(132, 13).Value = -12766.799
$ws.Cells.Item(132, 14).Value = -12514.4999
# Row 136: Metal with Mettle
$ws.Cells.Item(136, 8).Value = 2168.24
$ws.Cells.Item(136, 9).Value = 1860.3
$ws.Cells.Item(136, 11).Value = 5580.9
$ws.Cells.Item(136, 13).Value = -3030.9

$ws = $wb.Worksheets.Item("BSM")
# Row 86: Through Thick and Thin
$ws.Cells.Item(86, 8).Value = 57595.332
$ws.Cells.Item(86, 9).Value = 2308.3845
$ws.Cells.Item(86, 10).Value = 201341.4
$ws.Cells.Item(86, 11).Value = 2308.3845
$ws.Cells.Item(86, 12).Value = 201341.4
$ws.Cells.Item(86, 13).Value = -1185.3845
$ws.Cells.Item(86, 14).Value = -203587.4
# Row 89: Piercing Eyes Deserve Piercing Shafts (L)
$ws.Cells.Item(89, 8).Value = 57595.332
$ws.Cells.Item(89, 9).Value = 2308.3845
$ws.Cells.Item(89, 10).Value = 201341.4
$ws.Cells.Item(89, 11).Value = 11541.9225
$ws.Cells.Item(89, 12).Value = 1006707
$ws.Cells.Item(89, 13).Value = -5925.922500000001
$ws.Cells.Item(89, 14).Value = -1017939
# Row 94: High Steal
$ws.Cells.Item(94, 8).Value = 1960.091
$ws.Cells.Item(94, 9).Value = 1028
$ws.Cells.Item(94, 10).Value = 2736.8333
$ws.Cells.Item(94, 11).Value = 1028
$ws.Cells.Item(94, 12).Value = 2736.8333
$ws.Cells.Item(94, 13).Value = -577
$ws.Cells.Item(94, 14).Value = -3638.8333
# Row 99: Meddle in Metal
$ws.Cells.Item(99, 8).Value = 2109.889
$ws.Cells.Item(99, 9).Value = 1498.3334
$ws.Cells.Item(99, 11).Value = 1498.3334
$ws.Cells.Item(99, 13).Value = -0.3333999999999833
# Row 134: Ruthenium Supremium
$ws.Cells.Item(134, 8).Value = 2433.2
$ws.Cells.Item(134, 9).Value = 2177.577
$ws.Cells.Item(134, 10).Value = 2907.9285
$ws.Cells.Item(134, 11).Value = 6532.731000000001
$ws.Cells.Item(134, 12).Value = 8723.7855
$ws.Cells.Item(134, 13).Value = -3997.731000000001
$ws.Cells.Item(134, 14).Value = -13793.7855

$ws = $wb.Worksheets.Item("CRP")
# Row 31: Wall Not Found
$ws.Cells.Item(31, 8).Value = 2389.0312
$ws.Cells.Item(31, 9).Value = 1610.826
$ws.Cells.Item(31, 11).Value = 1610.826
$ws.Cells.Item(31, 13).Value = -1315.826
# Row 34: Armoires of the Rich and Famous
$ws.Cells.Item(34, 8).Value = 2389.0312
$ws.Cells.Item(34, 9).Value = 1610.826
$ws.Cells.Item(34, 11).Value = 1610.826
$ws.Cells.Item(34, 13).Value = -1408.826
# Row 58: You Do the Heavy Lifting
$ws.Cells.Item(58, 8).Value = 951690.1
$ws.Cells.Item(58, 9).Value = 1685477.1
$ws.Cells.Item(58, 11).Value = 1685477.1
$ws.Cells.Item(58, 13).Value = -1685274.1
# Row 107: Built to Last
$ws.Cells.Item(107, 8).Value = 409.6316
$ws.Cells.Item(107, 9).Value = 372.9091
$ws.Cells.Item(107, 10).Value = 460.125
$ws.Cells.Item(107, 11).Value = 372.9091
$ws.Cells.Item(107, 12).Value = 460.125
$ws.Cells.Item(107, 13).Value = 1547.0909
$ws.Cells.Item(107, 14).Value = -4300.125
# Row 122: Timber of Tenkonto
$ws.Cells.Item(122, 8).Value = 2202.0952
$ws.Cells.Item(122, 9).Value = 1862.4
$ws.Cells.Item(122, 11).Value = 5587.200000000001
$ws.Cells.Item(122, 13).Value = -3137.200000000001
# Row 127: In Rod We Trust
$ws.Cells.Item(127, 8).Value = 1833
$ws.Cells.Item(127, 10).Value = 1833
$ws.Cells.Item(127, 12).Value = 5499
$ws.Cells.Item(127, 14).Value = -15419
# Row 132: Hull Lotta Damage
$ws.Cells.Item(132, 8).Value = 502557.78
$ws.Cells.Item(132, 9).Value = 564690.0600000001
$ws.Cells.Item(132, 10).Value = 5499.3335
$ws.Cells.Item(132, 11).Value = 1694070.18
$ws.Cells.Item(132, 12).Value = 16498.0005
$ws.Cells.Item(132, 13).Value = -1691540.18
$ws.Cells.Item(132, 14).Value = -21558.0005
# Row 134: Wood You Be Quiet
$ws.Cells.Item(134, 8).Value = 3196.7896
$ws.Cells.Item(134, 9).Value = 2265.75
$ws.Cells.Item(134, 10).Value = 3873.9092
$ws.Cells.Item(134, 11).Value = 6797.25
$ws.Cells.Item(134, 12).Value = 11621.7276
$ws.Cells.Item(134, 13).Value = -4262.25
$ws.Cells.Item(134, 14).Value = -16691.7276
# Row 136: Turali Quality
$ws.Cells.Item(136, 8).Value = 951690.1
$ws.Cells.Item(136, 9).Value = 1685477.1
$ws.Cells.Item(136, 11).Value = 5056431.300000001
$ws.Cells.Item(136, 13).Value = -5053881.300000001

$ws = $wb.Worksheets.Item("CUL")
# Row 5: What a Sap
$ws.Cells.Item(5, 8).Value = 1741.6875
$ws.Cells.Item(5, 9).Value = 2213.6365
$ws.Cells.Item(5, 10).Value = 703.4
$ws.Cells.Item(5, 11).Value = 6640.9095
$ws.Cells.Item(5, 12).Value = 2110.2
$ws.Cells.Item(5, 13).Value = -6528.9095
$ws.Cells.Item(5, 14).Value = -2334.2
# Row 6: Meat-lover's Special
$ws.Cells.Item(6, 8).Value = 1858.25
$ws.Cells.Item(6, 9).Value = 238.66667
$ws.Cells.Item(6, 10).Value = 2830
$ws.Cells.Item(6, 11).Value = 716.00001
$ws.Cells.Item(6, 12).Value = 8490
$ws.Cells.Item(6, 13).Value = -603.00001
$ws.Cells.Item(6, 14).Value = -8716
# Row 132: More Mezcal
$ws.Cells.Item(132, 8).Value = 1630.25
$ws.Cells.Item(132, 9).Value = 1165.8889
$ws.Cells.Item(132, 11).Value = 10493.0001
$ws.Cells.Item(132, 13).Value = -7963.000099999999
# Row 135: Not-so-secret Ingredient
$ws.Cells.Item(135, 8).Value = 1741.6875
$ws.Cells.Item(135, 9).Value = 2213.6365
$ws.Cells.Item(135, 10).Value = 703.4
$ws.Cells.Item(135, 11).Value = 19922.7285
$ws.Cells.Item(135, 12).Value = 6330.599999999999
$ws.Cells.Item(135, 13).Value = -17387.7285
$ws.Cells.Item(135, 14).Value = -11400.6

$ws = $wb.Worksheets.Item("GSM")
# Row 132: On Board for Lar
$ws.Cells.Item(132, 8).Value = 2043.175
$ws.Cells.Item(132, 9).Value = 1374.9131
$ws.Cells.Item(132, 10).Value = 2947.2942
$ws.Cells.Item(132, 11).Value = 4124.7393
$ws.Cells.Item(132, 12).Value = 8841.882599999999
$ws.Cells.Item(132, 13).Value = -1594.7393
$ws.Cells.Item(132, 14).Value = -13901.8826

$ws = $wb.Worksheets.Item("LTW")
# Row 61: Spelling Me Softly
$ws.Cells.Item(61, 8).Value = 3200
$ws.Cells.Item(61, 10).Value = 4750
$ws.Cells.Item(61, 12).Value = 4750
$ws.Cells.Item(61, 14).Value = -5154
# Row 113: Peace in Rest
$ws.Cells.Item(113, 8).Value = 3200
$ws.Cells.Item(113, 10).Value = 4750
$ws.Cells.Item(113, 12).Value = 4750
$ws.Cells.Item(113, 14).Value = -9090
# Row 132: Tenets of Tanning
$ws.Cells.Item(132, 8).Value = 3241
$ws.Cells.Item(132, 9).Value = 2846.0952
$ws.Cells.Item(132, 10).Value = 4899.6
$ws.Cells.Item(132, 11).Value = 8538.285600000001
$ws.Cells.Item(132, 12).Value = 14698.8
$ws.Cells.Item(132, 13).Value = -6008.285600000001
$ws.Cells.Item(132, 14).Value = -19758.8
# Row 136: Respect for Br'aax
$ws.Cells.Item(136, 8).Value = 25252152
$ws.Cells.Item(136, 9).Value = 35715770
$ws.Cells.Item(136, 10).Value = 837051.25
$ws.Cells.Item(136, 11).Value = 107147310
$ws.Cells.Item(136, 12).Value = 2511153.75
$ws.Cells.Item(136, 13).Value = -107144760
$ws.Cells.Item(136, 14).Value = -2516253.75

$ws = $wb.Worksheets.Item("WVR")
# Row 126: A Polished Purchase
$ws.Cells.Item(126, 8).Value = 5235.905
$ws.Cells.Item(126, 9).Value = 5738.4707
$ws.Cells.Item(126, 10).Value = 3100
$ws.Cells.Item(126, 11).Value = 17215.4121
$ws.Cells.Item(126, 12).Value = 9300
$ws.Cells.Item(126, 13).Value = -14745.4121
$ws.Cells.Item(126, 14).Value = -14240
# Row 132: Comfy Cabins
$ws.Cells.Item(132, 8).Value = 1781.0465
$ws.Cells.Item(132, 9).Value = 1203.68
$ws.Cells.Item(132, 10).Value = 2582.9443
$ws.Cells.Item(132, 11).Value = 3611.04
$ws.Cells.Item(132, 12).Value = 7748.8329
$ws.Cells.Item(132, 13).Value = -1081.04
$ws.Cells.Item(132, 14).Value = -12808.8329
# Row 136: Weaving the Envelope
$ws.Cells.Item(136, 8).Value = 2579.4285
$ws.Cells.Item(136, 9).Value = 2801.7144
$ws.Cells.Item(136, 10).Value = 2357.1428
$ws.Cells.Item(136, 11).Value = 8405.143199999999
$ws.Cells.Item(136, 12).Value = 7071.428400000001
$ws.Cells.Item(136, 13).Value = -5855.143199999999
$ws.Cells.Item(136, 14).Value = -12171.4284
